$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet2 ("ThresholdConfig_PUT") restructuring
# ---------------------------------------------------------------------------
# Remove the orphan formatted cell at the bottom of the sheet.
$ws2.Range("A23").EntireRow.Delete()

# Remove the blank separator row after the first block (row 5).
$ws2.Range("A5").EntireRow.Delete()

# Remove the blank separator + placeholder rows before the 3rd block
# (originally rows 9 and 10, now rows 8 and 9 after the row5 delete above).
$ws2.Range("A8:A9").EntireRow.Delete()

# At this point the sheet looks like:
#  1-4  : Assert200 block (header/EndPoint/url/blank) - untouched
#  5-7  : Assert401 block (narrow header style)
#  8    : Assert400 label (narrow, single "header-ish" row)
#  9    : EndPoint   (compact, single cell)
#  10   : url        (compact, single cell)

# Make room: push the Assert400 block down so it starts on row 12
# (leaving rows 9-11 empty), matching the target layout.
$ws2.Range("A9:A11").EntireRow.Insert()

# Push things further down to create the new 4th block starting at row 18
# (leaving rows 15-17 empty).
$ws2.Range("A15:A17").EntireRow.Insert()

# ---------------------------------------------------------------------------
# Re-apply full-width header/EndPoint/url formatting by cloning the
# canonical block (rows 1-3) into every block position, then set the text.
# ---------------------------------------------------------------------------

function Set-Block($headerRow, $epRow, $urlRow, $headerText)
{
    $ws2.Range("A1:E1").Copy()
    $ws2.Range("A" + $headerRow + ":E" + $headerRow).PasteSpecial(-4122)
    $ws2.Range("A2:E2").Copy()
    $ws2.Range("A" + $epRow + ":E" + $epRow).PasteSpecial(-4122)
    $ws2.Range("A3:E3").Copy()
    $ws2.Range("A" + $urlRow + ":E" + $urlRow).PasteSpecial(-4122)

    $ws2.Range("A" + $headerRow).Value = $headerText
    $ws2.Range("A" + $epRow).Value = "EndPoint"
    $ws2.Range("A" + $urlRow).Value = "/configuration/system/instrumentThreshold"
}

# Block 2: Assert401 (rows 6,7,8)
Set-Block 6 7 8 "Assert401"

# Block 3: Assert400 (rows 12,13,14)
Set-Block 12 13 14 "Assert400"

# Block 4: Assert409 (rows 18,19,20) - brand new block
Set-Block 18 19 20 "Assert409"

# ---------------------------------------------------------------------------
# Selections (the author left sheet1 with A11 selected, then returned to
# sheet2 with A2 selected and sheet2 remaining the active/visible tab).
# ---------------------------------------------------------------------------
$ws1.Range("A11").Select()
$ws2.Activate()
$ws2.Range("A2").Select()
